$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 710, shifting existing rows 710-743 down to 711-744.
$ws.Rows("710:710").Insert()

# Populate the newly inserted row 710 with the new data record.
$ws.Range("A710").Value = 10
$ws.Range("B710").Value = "Vega Modelo de Temuco"
$ws.Range("C710").Value = "La Araucanía"
$ws.Range("D710").Value = 45147
$ws.Range("E710").Value = 9
$ws.Range("F710").Value = 100112023
$ws.Range("G710").Value = "Brócoli"
$ws.Range("H710").Value = "Sin especificar"
$ws.Range("I710").Value = "Primera"
$ws.Range("J710").Value = 1260
$ws.Range("K710").Value = 1000
$ws.Range("L710").Value = 1000
$ws.Range("M710").Value = 1000
$ws.Range("N710").Value = "$/unidad"
$ws.Range("O710").Value = "Región del Maule"
$ws.Range("P710").Value = 1000
$ws.Range("Q710").Value = 1
$ws.Range("R710").Value = "Hortaliza"
